# Daily update at 8 AM UTC
# Appends the next day's row (row 75) to the "Wins Over Time" tracking
# sheet, and moves the "last row" number formatting down from the old
# last row (74) to the new one (75).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 74 was previously the last row (date-only format). It is no longer
# the last row, so it goes back to the regular "date + time" format used
# by every other data row.
$ws.Cells.Item(74, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New data row for 2025-06-07 (serial 45815).
$ws.Cells.Item(75, 1).Value = 45815
$ws.Cells.Item(75, 2).Value = 321
$ws.Cells.Item(75, 3).Value = 318
$ws.Cells.Item(75, 4).Value = 324

# The new last row gets the date-only format.
$ws.Cells.Item(75, 1).NumberFormat = "YYYY-MM-DD"
